$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="37.141.10"; E="  -0.25%  " },
    @{ Row=3; D="2.073.12"; E="  -0.74%  " },
    @{ Row=4; E="  +0.11%  " },
    @{ Row=5; D="253.40"; E="  +1.14%  " },
    @{ Row=6; D="0.675"; E="  +1.93%  " },
    @{ Row=7; D="59.99"; E="  +10.28%  " },
    @{ Row=8; E="  -0.07%  " },
    @{ Row=9; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="61.56"; E="  -0.21%  " },
    @{ Row=10; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.389"; E="  +4.23%  " },
    @{ Row=11; D="0.0800"; E="  +7.81%  " },
    @{ Row=12; E="  +2.48%  " },
    @{ Row=13; E="  +8.01%  " },
    @{ Row=14; D="2.379.18"; E="  -0.62%  " },
    @{ Row=15; D="0.816"; E="  -1.81%  " },
    @{ Row=16; D="5.65"; E="  +9.68%  " },
    @{ Row=17; D="2.073.83"; E="  -0.72%  " },
    @{ Row=18; D="37.150.57"; E="  -0.13%  " },
    @{ Row=19; D="16.68"; E="  +14.47%  " },
    @{ Row=20; D="75.06"; E="  +3.22%  " },
    @{ Row=21; D="0.0₃0933"; E="  +10.46%  " },
    @{ Row=22; D="5.50"; E="  +5.91%  " },
    @{ Row=23; D="239.61"; E="  -0.33%  " },
    @{ Row=24; E="  -0.05%  " },
    @{ Row=25; D="2.41"; E="  -2.30%  " },
    @{ Row=26; E="  +14.69%  " },
    @{ Row=27; D="169.51"; E="  -1.56%  " },
    @{ Row=28; D="9.38"; E="  +1.62%  " },
    @{ Row=29; D="20.42"; E="  -1.30%  " },
    @{ Row=30; E="  +3.26%  " },
    @{ Row=31; E="  +5.01%  " },
    @{ Row=32; D="4.80"; E="  +6.34%  " },
    @{ Row=33; D="0.0625"; E="  +1.61%  " },
    @{ Row=34; D="4.55"; E="  +10.03%  " },
    @{ Row=35; D="0.0915"; E="  +0.06%  " },
    @{ Row=36; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.00"; E="  +0.09%  " },
    @{ Row=37; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="2.32"; E="  +3.92%  " },
    @{ Row=38; E="  +26.08%  " },
    @{ Row=39; E="  -3.91%  " },
    @{ Row=40; E="  +2.82%  " },
    @{ Row=41; D="18.11"; E="  -0.59%  " },
    @{ Row=42; E="  +0.82%  " },
    @{ Row=43; E="  +0.37%  " },
    @{ Row=44; D="98.94"; E="  +0.57%  " },
    @{ Row=45; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.86"; E="  +2.38%  " },
    @{ Row=46; B="FTXToken"; C="https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; D="4.33"; E="  -0.06%  " },
    @{ Row=47; D="4.59"; E="  +13.97%  " },
    @{ Row=48; D="2.52"; E="  +7.76%  " },
    @{ Row=49; D="1.305.67"; E="  -1.15%  " },
    @{ Row=50; D="2.93"; E="  -0.29%  " },
    @{ Row=51; D="6.97"; E="  -0.44%  " },
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
